$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two student rows (19 and 20) - roster shrank from 19 to 17 students
$ws.Rows("19:20").Delete()

# New roster data: StudentID, Name, Score for rows 2-18
$data = @(
    @(2201804568, "符群舒", 78.75),
    @(2201804800, "孔兴",   75),
    @(2201804904, "桂云伯", 75),
    @(2201804856, "苗淑媛", 73.75),
    @(2201804433, "黄寒",   72.5),
    @(2201804854, "阮德以", 72.5),
    @(2201804914, "幸辉",   72.5),
    @(2201804488, "王妍",   71.25),
    @(2201804505, "姜梦",   71.25),
    @(2201804804, "牛裕环", 71.25),
    @(2201804978, "盛腾鹏", 71.25),
    @(2201804443, "黄菁航", 70),
    @(2201804465, "向悦",   70),
    @(2201804645, "向天纨", 70),
    @(2201804694, "桓军伦", 70),
    @(2201804727, "牧兴",   70),
    @(2201804921, "阮丹毓", 70)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
